$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.032369287002124
$ws.Cells.Item(2, 4).Value = 1.038877566845704
$ws.Cells.Item(2, 5).Value = 1.03600998851591
$ws.Cells.Item(2, 6).Value = 1.046290008889653
$ws.Cells.Item(2, 9).Value = 1.034252192964906
$ws.Cells.Item(2, 10).Value = 1.037499885917508
$ws.Cells.Item(2, 11).Value = 1.04166434889497
$ws.Cells.Item(2, 12).Value = 1.03880495523503
$ws.Cells.Item(2, 13).Value = 1.049055859801566
$ws.Cells.Item(2, 14).Value = 1.038973254758411

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.033467922030577
$ws.Cells.Item(3, 4).Value = 1.039719290054452
$ws.Cells.Item(3, 5).Value = 1.037056905091805
$ws.Cells.Item(3, 6).Value = 1.047333409820946
$ws.Cells.Item(3, 9).Value = 1.034474207408196
$ws.Cells.Item(3, 10).Value = 1.038240029100217
$ws.Cells.Item(3, 11).Value = 1.042316092220624
$ws.Cells.Item(3, 12).Value = 1.03966075458725
$ws.Cells.Item(3, 13).Value = 1.049910270949507
$ws.Cells.Item(3, 14).Value = 1.039714449029335

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.034178759084533
$ws.Cells.Item(4, 4).Value = 1.040263524802836
$ws.Cells.Item(4, 5).Value = 1.037734635416843
$ws.Cells.Item(4, 6).Value = 1.048008636136325
$ws.Cells.Item(4, 9).Value = 1.034615931612829
$ws.Cells.Item(4, 10).Value = 1.038718358032192
$ws.Cells.Item(4, 11).Value = 1.042736755084068
$ws.Cells.Item(4, 12).Value = 1.040214221299974
$ws.Cells.Item(4, 13).Value = 1.050462606796468
$ws.Cells.Item(4, 14).Value = 1.040193457243258

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.034477582810916
$ws.Cells.Item(5, 4).Value = 1.040492221021307
$ws.Cells.Item(5, 5).Value = 1.03801962637527
$ws.Cells.Item(5, 6).Value = 1.048292519391357
$ws.Cells.Item(5, 9).Value = 1.034675049199905
$ws.Cells.Item(5, 10).Value = 1.038919305512804
$ws.Cells.Item(5, 11).Value = 1.042913348101665
$ws.Cells.Item(5, 12).Value = 1.040446828657828
$ws.Cells.Item(5, 13).Value = 1.050694682894506
$ws.Cells.Item(5, 14).Value = 1.040394690092342

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.03452775593266
$ws.Cells.Item(6, 4).Value = 1.040530614243659
$ws.Cells.Item(6, 5).Value = 1.038067481900723
$ws.Cells.Item(6, 6).Value = 1.048340185697633
$ws.Cells.Item(6, 9).Value = 1.034684948128706
$ws.Cells.Item(6, 10).Value = 1.038953037152156
$ws.Cells.Item(6, 11).Value = 1.042942983966042
$ws.Cells.Item(6, 12).Value = 1.040485880349682
$ws.Cells.Item(6, 13).Value = 1.050733642106165
$ws.Cells.Item(6, 14).Value = 1.040428469634491

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.034182752028892
$ws.Cells.Item(7, 4).Value = 1.040266581046593
$ws.Cells.Item(7, 5).Value = 1.037738443192588
$ws.Cells.Item(7, 6).Value = 1.048012429327269
$ws.Cells.Item(7, 9).Value = 1.034616723365601
$ws.Cells.Item(7, 10).Value = 1.038721043660237
$ws.Cells.Item(7, 11).Value = 1.042739115725566
$ws.Cells.Item(7, 12).Value = 1.040217329686101
$ws.Cells.Item(7, 13).Value = 1.050465708303506
$ws.Cells.Item(7, 14).Value = 1.040196146685204

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.032740587235901
$ws.Cells.Item(8, 4).Value = 1.039162117154203
$ws.Cells.Item(8, 5).Value = 1.036363735687647
$ws.Cells.Item(8, 6).Value = 1.046642615285807
$ws.Cells.Item(8, 9).Value = 1.034327624070508
$ws.Cells.Item(8, 10).Value = 1.037750143605434
$ws.Cells.Item(8, 11).Value = 1.04188482768123
$ws.Cells.Item(8, 12).Value = 1.039094237252895
$ws.Cells.Item(8, 13).Value = 1.049344720737283
$ws.Cells.Item(8, 14).Value = 1.039223867840959

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.030198873914291
$ws.Cells.Item(9, 4).Value = 1.037212736649512
$ws.Cells.Item(9, 5).Value = 1.033943660952538
$ws.Cells.Item(9, 6).Value = 1.044229415343362
$ws.Cells.Item(9, 9).Value = 1.033803393490287
$ws.Cells.Item(9, 10).Value = 1.036034749249194
$ws.Cells.Item(9, 11).Value = 1.0403713633929
$ws.Cells.Item(9, 12).Value = 1.037112961396655
$ws.Cells.Item(9, 13).Value = 1.047365379354678
$ws.Cells.Item(9, 14).Value = 1.037506037427974

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.028504065639006
$ws.Cells.Item(10, 4).Value = 1.035911033489329
$ws.Cells.Item(10, 5).Value = 1.032331840843971
$ws.Cells.Item(10, 6).Value = 1.042621012396929
$ws.Cells.Item(10, 9).Value = 1.033443968942215
$ws.Cells.Item(10, 10).Value = 1.034888092954161
$ws.Cells.Item(10, 11).Value = 1.039356949960692
$ws.Cells.Item(10, 12).Value = 1.035790596717775
$ws.Cells.Item(10, 13).Value = 1.046043120477531
$ws.Cells.Item(10, 14).Value = 1.036357752749479

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.027770106039837
$ws.Cells.Item(11, 4).Value = 1.035346882519169
$ws.Cells.Item(11, 5).Value = 1.031634272981608
$ws.Cells.Item(11, 6).Value = 1.041924649591113
$ws.Cells.Item(11, 9).Value = 1.033285979444769
$ws.Cells.Item(11, 10).Value = 1.034390850357543
$ws.Cells.Item(11, 11).Value = 1.038916408979606
$ws.Cells.Item(11, 12).Value = 1.03521763682739
$ws.Cells.Item(11, 13).Value = 1.045469927554834
$ws.Cells.Item(11, 14).Value = 1.03585980401134

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.027497465301857
$ws.Cells.Item(12, 4).Value = 1.035137255984262
$ws.Cells.Item(12, 5).Value = 1.031375219087168
$ws.Cells.Item(12, 6).Value = 1.041666002209173
$ws.Cells.Item(12, 9).Value = 1.033226941213382
$ws.Cells.Item(12, 10).Value = 1.034206041938619
$ws.Cells.Item(12, 11).Value = 1.038752578101805
$ws.Cells.Item(12, 12).Value = 1.03500475853005
$ws.Cells.Item(12, 13).Value = 1.04525692096057
$ws.Cells.Item(12, 14).Value = 1.035674733143263

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.0275559483811
$ws.Cells.Item(13, 4).Value = 1.035182225013543
$ws.Cells.Item(13, 5).Value = 1.031430784594891
$ws.Cells.Item(13, 6).Value = 1.04172148237592
$ws.Cells.Item(13, 9).Value = 1.033239621127333
$ws.Cells.Item(13, 10).Value = 1.034245688971179
$ws.Cells.Item(13, 11).Value = 1.038787729176393
$ws.Cells.Item(13, 12).Value = 1.035050424163161
$ws.Cells.Item(13, 13).Value = 1.045302616000082
$ws.Cells.Item(13, 14).Value = 1.035714436479157

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.027747569792484
$ws.Cells.Item(14, 4).Value = 1.035329556268975
$ws.Cells.Item(14, 5).Value = 1.031612858409012
$ws.Cells.Item(14, 6).Value = 1.041903269458482
$ws.Cells.Item(14, 9).Value = 1.033281106544934
$ws.Cells.Item(14, 10).Value = 1.034375576289291
$ws.Cells.Item(14, 11).Value = 1.038902870632938
$ws.Cells.Item(14, 12).Value = 1.035200041374553
$ws.Cells.Item(14, 13).Value = 1.045452322347491
$ws.Cells.Item(14, 14).Value = 1.035844508252158

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.027865631990075
$ws.Cells.Item(15, 4).Value = 1.035420321884141
$ws.Cells.Item(15, 5).Value = 1.031725047222696
$ws.Cells.Item(15, 6).Value = 1.042015276165988
$ws.Cells.Item(15, 9).Value = 1.033306620188232
$ws.Cells.Item(15, 10).Value = 1.034455589518206
$ws.Cells.Item(15, 11).Value = 1.038973787330814
$ws.Cells.Item(15, 12).Value = 1.035292218115189
$ws.Cells.Item(15, 13).Value = 1.045544548481642
$ws.Cells.Item(15, 14).Value = 1.035924635109036

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.028552774005462
$ws.Cells.Item(16, 4).Value = 1.035948463703023
$ws.Cells.Item(16, 5).Value = 1.032378143699245
$ws.Cells.Item(16, 6).Value = 1.042667229504471
$ws.Cells.Item(16, 9).Value = 1.033454404549829
$ws.Cells.Item(16, 10).Value = 1.034921077843847
$ws.Cells.Item(16, 11).Value = 1.039386159963473
$ws.Cells.Item(16, 12).Value = 1.035828614414933
$ws.Cells.Item(16, 13).Value = 1.046081147778595
$ws.Cells.Item(16, 14).Value = 1.036390784481492

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.028983773532892
$ws.Cells.Item(17, 4).Value = 1.036279618014986
$ws.Cells.Item(17, 5).Value = 1.032787910237704
$ws.Cells.Item(17, 6).Value = 1.043076205575771
$ws.Cells.Item(17, 9).Value = 1.033546475006974
$ws.Cells.Item(17, 10).Value = 1.035212869973774
$ws.Cells.Item(17, 11).Value = 1.039644484147914
$ws.Cells.Item(17, 12).Value = 1.036164983102309
$ws.Cells.Item(17, 13).Value = 1.046417569255909
$ws.Cells.Item(17, 14).Value = 1.036682990989713

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.02923515902666
$ws.Cells.Item(18, 4).Value = 1.036472725998784
$ws.Cells.Item(18, 5).Value = 1.033026954970416
$ws.Cells.Item(18, 6).Value = 1.043314762737344
$ws.Cells.Item(18, 9).Value = 1.033599950791456
$ws.Cells.Item(18, 10).Value = 1.035382996695418
$ws.Cells.Item(18, 11).Value = 1.039795035495594
$ws.Cells.Item(18, 12).Value = 1.036361145852519
$ws.Cells.Item(18, 13).Value = 1.046613735882777
$ws.Cells.Item(18, 14).Value = 1.036853359310815

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.029320873446236
$ws.Cells.Item(19, 4).Value = 1.036538562562999
$ws.Cells.Item(19, 5).Value = 1.033108468967417
$ws.Cells.Item(19, 6).Value = 1.043396105984706
$ws.Cells.Item(19, 9).Value = 1.033618146085191
$ws.Cells.Item(19, 10).Value = 1.035440993528921
$ws.Cells.Item(19, 11).Value = 1.039846348451146
$ws.Cells.Item(19, 12).Value = 1.036428026283424
$ws.Cells.Item(19, 13).Value = 1.046680613065562
$ws.Cells.Item(19, 14).Value = 1.036911438506475

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.028937532324973
$ws.Cells.Item(20, 4).Value = 1.036244093330095
$ws.Cells.Item(20, 5).Value = 1.032743942564446
$ws.Cells.Item(20, 6).Value = 1.043032325452568
$ws.Cells.Item(20, 9).Value = 1.033536620233842
$ws.Cells.Item(20, 10).Value = 1.035181570753478
$ws.Cells.Item(20, 11).Value = 1.039616781314291
$ws.Cells.Item(20, 12).Value = 1.036128897590637
$ws.Cells.Item(20, 13).Value = 1.04638148086318
$ws.Cells.Item(20, 14).Value = 1.036651647320934

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.027691142497581
$ws.Cells.Item(21, 4).Value = 1.035286172976483
$ws.Cells.Item(21, 5).Value = 1.03155924072461
$ws.Cells.Item(21, 6).Value = 1.041849737342671
$ws.Cells.Item(21, 9).Value = 1.033268899891899
$ws.Cells.Item(21, 10).Value = 1.034337330760697
$ws.Cells.Item(21, 11).Value = 1.038868969708989
$ws.Cells.Item(21, 12).Value = 1.035155984368997
$ws.Cells.Item(21, 13).Value = 1.045408240244546
$ws.Cells.Item(21, 14).Value = 1.035806208410527

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.026907397793023
$ws.Cells.Item(22, 4).Value = 1.034683452238234
$ws.Cells.Item(22, 5).Value = 1.0308146831493
$ws.Cells.Item(22, 6).Value = 1.041106271042256
$ws.Cells.Item(22, 9).Value = 1.033098525995495
$ws.Cells.Item(22, 10).Value = 1.033805884811696
$ws.Cells.Item(22, 11).Value = 1.038397666441209
$ws.Cells.Item(22, 12).Value = 1.034543954048879
$ws.Cells.Item(22, 13).Value = 1.044795762667313
$ws.Cells.Item(22, 14).Value = 1.035274007747321

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.027322884542855
$ws.Cells.Item(23, 4).Value = 1.035003007405924
$ws.Cells.Item(23, 5).Value = 1.031209357799685
$ws.Cells.Item(23, 6).Value = 1.041500389667772
$ws.Cells.Item(23, 9).Value = 1.033189038464987
$ws.Cells.Item(23, 10).Value = 1.034087675004844
$ws.Cells.Item(23, 11).Value = 1.03864761982009
$ws.Cells.Item(23, 12).Value = 1.034868433447851
$ws.Cells.Item(23, 13).Value = 1.045120501997254
$ws.Cells.Item(23, 14).Value = 1.035556198114865

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.028958426764601
$ws.Cells.Item(24, 4).Value = 1.036260145556468
$ws.Cells.Item(24, 5).Value = 1.032763809555542
$ws.Cells.Item(24, 6).Value = 1.043052152965477
$ws.Cells.Item(24, 9).Value = 1.033541073884045
$ws.Cells.Item(24, 10).Value = 1.035195713742105
$ws.Cells.Item(24, 11).Value = 1.039629299417745
$ws.Cells.Item(24, 12).Value = 1.036145203189299
$ws.Cells.Item(24, 13).Value = 1.046397787846848
$ws.Cells.Item(24, 14).Value = 1.036665810394227

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.030856023398488
$ws.Cells.Item(25, 4).Value = 1.037717072498249
$ws.Cells.Item(25, 5).Value = 1.034569032231015
$ws.Cells.Item(25, 6).Value = 1.04485321551902
$ws.Cells.Item(25, 9).Value = 1.033940672189582
$ws.Cells.Item(25, 10).Value = 1.036478758801539
$ws.Cells.Item(25, 11).Value = 1.040763589131045
$ws.Cells.Item(25, 12).Value = 1.037625435128279
$ws.Cells.Item(25, 13).Value = 1.047877562312412
$ws.Cells.Item(25, 14).Value = 1.037950677524812
